# "Generate Report for Handback" — the CI localization-status report is
# regenerated after a handback: file "a.md" has come back in sync for both
# the zh-cn and de-de targets, so the Status cells flip from
# "Ready for handoff" to "Handed back: in sync with en-US", and the
# per-language detail sheets get their Latest Target File / Latest Handback
# File / Latest Handback DateTime columns populated (with a hyperlink on the
# target-file cell, matching the existing source-file hyperlink style).

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$srcUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8580e004da770ffceef172cdbe1908c444750cf0/e2e/a.md"

# zh-cn detail sheet
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus

$zh.Range("I2").Value = "a.md"
$zh.Hyperlinks.Add($zh.Range("I2"), $srcUrl, "", "", "a.md")
$zh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-19 20:42:26"

$zh.Range("I3").Value = "a.md"
$zh.Hyperlinks.Add($zh.Range("I3"), $srcUrl, "", "", "a.md")
$zh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-19 20:42:26"

$zh.Columns.Item(3).ColumnWidth = 29.17
$zh.Columns.Item(10).ColumnWidth = 40

# de-de detail sheet
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus

$de.Range("I2").Value = "a.md"
$de.Hyperlinks.Add($de.Range("I2"), $srcUrl, "", "", "a.md")
$de.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$de.Range("K2").Value = "2016-08-19 20:42:32"

$de.Range("I3").Value = "a.md"
$de.Hyperlinks.Add($de.Range("I3"), $srcUrl, "", "", "a.md")
$de.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$de.Range("K3").Value = "2016-08-19 20:42:32"

$de.Columns.Item(3).ColumnWidth = 29.17
$de.Columns.Item(10).ColumnWidth = 40

# Overview Status columns (E/F) auto-widen along with the detail-sheet
# Status column now that the text is longer.
$overview.Columns.Item(5).ColumnWidth = 29.17
$overview.Columns.Item(6).ColumnWidth = 29.17
